$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C ("username"), shifting old C..I to D..J
$ws.Columns("C:C").Insert()

# New header for the inserted column (style already carried over by Insert)
$ws.Cells.Item(1, 3).Value = "username"

# Make sure date/time text columns (D = tanggal tugas, E = terakhir submit)
# are stored as plain text instead of being auto-converted into Excel
# date/time serial numbers.
$ws.Range("D2:E7").NumberFormat = "@"

# Update existing row 2 with corrected data
$ws.Cells.Item(2, 1).Value = 507549293
$ws.Cells.Item(2, 2).Value = "Kharisma Muzaki"
$ws.Cells.Item(2, 3).Value = "muzaki_gh"
$ws.Cells.Item(2, 4).Value = "2020-02-04"
$ws.Cells.Item(2, 5).Value = "15:41:35"
$ws.Cells.Item(2, 6).Value = 5
$ws.Cells.Item(2, 7).Value = 6
$ws.Cells.Item(2, 8).Value = 1
$ws.Cells.Item(2, 9).Value = 11
$ws.Cells.Item(2, 10).Value = 12

# Row 3
$ws.Cells.Item(3, 1).Value = 507549293
$ws.Cells.Item(3, 2).Value = "Kharisma Muzaki"
$ws.Cells.Item(3, 3).Value = "muzaki_gh"
$ws.Cells.Item(3, 4).Value = "2020-02-05"
$ws.Cells.Item(3, 5).Value = "09:53:55"
$ws.Cells.Item(3, 6).Value = 4
$ws.Cells.Item(3, 7).Value = 10
$ws.Cells.Item(3, 8).Value = 1
$ws.Cells.Item(3, 9).Value = 14
$ws.Cells.Item(3, 10).Value = 15

# Row 4
$ws.Cells.Item(4, 1).Value = 123123087
$ws.Cells.Item(4, 2).Value = "Harris Setyawan"
$ws.Cells.Item(4, 3).Value = "harris"
$ws.Cells.Item(4, 4).Value = "2020-02-04"
$ws.Cells.Item(4, 5).Value = "15:41:35"
$ws.Cells.Item(4, 6).Value = 10
$ws.Cells.Item(4, 7).Value = 6
$ws.Cells.Item(4, 8).Value = 5
$ws.Cells.Item(4, 9).Value = 16
$ws.Cells.Item(4, 10).Value = 21

# Row 5
$ws.Cells.Item(5, 1).Value = 3122331
$ws.Cells.Item(5, 2).Value = "Riko Alfianto"
$ws.Cells.Item(5, 3).Value = "riko"
$ws.Cells.Item(5, 4).Value = "2020-02-04"
$ws.Cells.Item(5, 5).Value = "15:41:35"
$ws.Cells.Item(5, 6).Value = 13
$ws.Cells.Item(5, 7).Value = 7
$ws.Cells.Item(5, 8).Value = 1
$ws.Cells.Item(5, 9).Value = 20
$ws.Cells.Item(5, 10).Value = 21

# Row 6
$ws.Cells.Item(6, 1).Value = 507123087
$ws.Cells.Item(6, 2).Value = "Anada Badu"
$ws.Cells.Item(6, 3).Value = "anada"
$ws.Cells.Item(6, 4).Value = "2020-02-04"
$ws.Cells.Item(6, 5).Value = "15:41:35"
$ws.Cells.Item(6, 6).Value = 8
$ws.Cells.Item(6, 7).Value = 5
$ws.Cells.Item(6, 8).Value = 6
$ws.Cells.Item(6, 9).Value = 13
$ws.Cells.Item(6, 10).Value = 19

# Row 7
$ws.Cells.Item(7, 1).Value = 401123087
$ws.Cells.Item(7, 2).Value = "Mozaze"
$ws.Cells.Item(7, 3).Value = "moza"
$ws.Cells.Item(7, 4).Value = "2020-02-04"
$ws.Cells.Item(7, 5).Value = "15:41:35"
$ws.Cells.Item(7, 6).Value = 12
$ws.Cells.Item(7, 7).Value = 1
$ws.Cells.Item(7, 8).Value = 2
$ws.Cells.Item(7, 9).Value = 13
$ws.Cells.Item(7, 10).Value = 15

# The original cells had no special number format (default "Normal" style);
# restore that now that the text values are safely stored, so no stray
# "@" text-format style is left applied to these cells.
$ws.Range("D2:E7").Style = "Normal"

Write-Host "done"
